# Commit: "Add files via upload"
# The authoritative (non-environment-noise) change in the diff is the
# renaming of the "northeast" worksheet to "northeastern" (see
# xl/workbook.xml <sheets> list). Everything else in the upstream diff
# (revisionPtr documentId GUID, bookViews window geometry, per-sheet
# x14ac:dyDescent / defaultRowHeight tweaks, and the disappearance of the
# <selection> elements on the "north" and "northeast" sheets) is the kind
# of incidental re-save noise Excel produces when a file is opened/saved
# on a different machine/session, not a deliberate edit, so it is not
# reproduced here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("northeast")
$ws.Name = "northeastern"
